$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
